$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Narrow column B (title) from 51 to 45 characters (45 in saved OOXML width units)
$ws.Columns.Item(2).ColumnWidth = 44.17

# Row 2
$ws.Range("A2").Value = '2025-12-19 12:37:51'
$ws.Range("B2").Value = 'EC×AIプロダクト/業務改善リード'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5450024'
$ws.Range("G2").Value = 338
$ws.Range("H2").Value = '🔥AI,Ai ◇業務改善'

# Row 3
$ws.Range("A3").Value = '2025-12-19 12:37:51'
$ws.Range("B3").Value = '製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5439165'
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-12-19 12:37:51'
$ws.Range("B4").Value = '【急募】AWSスクレイピングツールの開発を依頼したいです!'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5457255'
$ws.Range("G4").Value = 165
$ws.Range("H4").Value = '◆ツール,開発'

# Row 5
$ws.Range("A5").Value = '2025-12-19 12:37:51'
$ws.Range("B5").Value = '【フリーランス募集】Webサービス・業務システム開発エンジニア'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5457382'
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = '◆開発,システム開発'

# Row 6
$ws.Range("A6").Value = '2025-12-19 12:37:51'
$ws.Range("B6").Value = '【急募】飲食店予約サイトの制作と将来的なアプリ化(アプリ化の際は別契約)'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5457089'
$ws.Range("G6").Value = 70
$ws.Range("H6").Value = '◇アプリ'

# Row 7
$ws.Range("A7").Value = '2025-12-19 12:37:51'
$ws.Range("B7").Value = '【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5457026'
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = '◆ツール'

# Row 8
$ws.Range("A8").Value = '2025-12-19 12:37:51'
$ws.Range("B8").Value = '【急募】PHPによる申請サイト構築支援!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5457023'
$ws.Range("G8").Value = 58
$ws.Range("H8").Value = '◇サイト ○PHP'

# Row 9
$ws.Range("A9").Value = '2025-12-19 12:37:51'
$ws.Range("B9").Value = '【急募】Kintoneでの請求書自動発行システム構築依頼'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5457134'
$ws.Range("G9").Value = 28

# Row 10
$ws.Range("A10").Value = '2025-12-19 12:37:51'
$ws.Range("B10").Value = '回路設計者募集|UVA浄化装置 (マイコン不使用/555タイマー制御) ※成果物全帰属'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5457451'
$ws.Range("G10").Value = 18

# Rebuild hyperlinks: the scraper re-sorted/shifted rows, so every F-column
# hyperlink target needs to line up with its (possibly new) row. This emulation
# only supports clearing ALL hyperlinks on the sheet at once, so drop them all
# and re-add the full set with the correct targets.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5450024')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5439165')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5457255')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5457382')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5457089')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5457026')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5457023')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5457134')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5457451')

# Restore the Hyperlink cell style (Hyperlinks.Add resets it on each call)
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"

"done"